$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Test Steps" sheet: add a new "Result1" column (H) that records PASS for
# every real data row (the blank separator row 15 is left untouched).
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Steps")

$wsSteps.Range("H1").Value = "Result1"
$wsSteps.Range("H1").Interior.ColorIndex = 55

$dataRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,22,23,24,25,26)
foreach ($r in $dataRows) {
    $cell = $wsSteps.Cells.Item($r, 8)
    $cell.Value = "PASS"
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet selection / active-tab bookkeeping: "Test Cases" becomes the active
# sheet (was "Test Steps"), with the last selection sitting on B15.
# "Test Steps" keeps its own last selection on H10.
# ---------------------------------------------------------------------------
$wsSteps.Range("H10").Select() | Out-Null

$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Activate()
$wsCases.Range("B15").Select() | Out-Null
